$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51: new paper entry
$ws.Cells.Item(51, 1).Value = "室内场景下实时地三维语义地图构建"
$ws.Cells.Item(51, 2).Value = "单吉超，李秀智"
$ws.Cells.Item(51, 3).Value = 2019
$ws.Cells.Item(51, 4).Value = "仪器仪表学报"
$ws.Cells.Item(51, 5).Value = "基于深度残差网络( DＲN) 的像素级图像语义分割和三维同时定位与建图( SLAM) 相融合的三维语义地图构建方法。首先，采用一种联合中值滤波算法进行深度图像的修复，使用改进的迭代最近点( ICP) 算法得到相机估计`n位姿以及基于随机蕨类的闭环检测构建出三维环境地图; 其次，采用优化的深度残差网络对输入的图像实现较精准的像素语义级别的预测与分割; 最后，采用贝叶斯更新方法，渐进式的将图像分割获取的语义分类标签迁移到重建的室内三维模型中，获得完整的三维语义地图。实验表明，所设计的方法可以在实际的、复杂环境下实时地构建语义地图。"

# Row 52: new paper entry
$ws.Cells.Item(52, 1).Value = "实时的移动机器人语义地图构建系统"
$ws.Cells.Item(52, 2).Value = "李秀智，李尚宇"
$ws.Cells.Item(52, 3).Value = 2017
$ws.Cells.Item(52, 4).Value = "仪器仪表学报"
